$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp footer (row 1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 09:22"

# 2. Row 26 - Israel
$ws.Range("B26").Value = 15398
$ws.Range("C26").Value = 100
$ws.Range("D26").Value = 6602
$ws.Range("E26").Value = 8597
$ws.Range("F26").Value = 132

# 3. Row 33 - Polonia
$ws.Range("D33").Value = 2265
$ws.Range("E33").Value = 8484

# 4. Row 44 - Chequia
$ws.Range("D44").Value = 2471
$ws.Range("E44").Value = 4662
$ws.Range("F44").Value = 79
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 219

# 5. Row 71 - Armenia
$ws.Range("B71").Value = 1746
$ws.Range("C71").Value = 69
$ws.Range("D71").Value = 833
$ws.Range("E71").Value = 885

# 6 & 7. Rows 92/93 - Letonia overtakes Republica de Chipre in the ranking,
# so the two countries swap row order (row 92 becomes Letonia, row 93
# becomes Republica de Chipre) while their case counts are updated.
$ws.Range("A92").Value = "Letonia"
$ws.Range("B92").Value = 812
$ws.Range("C92").Value = 8
$ws.Range("D92").Value = 267
$ws.Range("E92").Value = 533
$ws.Range("F92").Value = 6
$ws.Range("H92").Value = 12

$ws.Range("A93").Value = "Republica de Chipre"
$ws.Range("B93").Value = 810
$ws.Range("D93").Value = 148
$ws.Range("E93").Value = 648
$ws.Range("F93").Value = 15
$ws.Range("H93").Value = 14

# 8. Row 107 - Sri Lanka
$ws.Range("B107").Value = 462
$ws.Range("C107").Value = 10
$ws.Range("E107").Value = 337

# 9. Row 167 - Nepal
$ws.Range("B167").Value = 51
$ws.Range("C167").Value = 2
$ws.Range("E167").Value = 39
